# Reorders the three requirement lines in the "Requisitos" bulleted
# paragraph so that "LOB1012 -  Estatística  (Requisito fraco)" moves
# from the first line to the last line, i.e.:
#   LOB1012, LOQ4095, LOQ4098  ->  LOQ4095, LOQ4098, LOB1012
#
# The paragraph holds three runs, each "<text><w:br/>", one per line.
# A plain in-place text swap causes Word to merge adjacent runs that
# get edited together, which would collapse the three distinct <w:r>
# elements into one. To preserve the original run-per-line structure,
# we build the new line order in three brand-new paragraphs appended
# at the end of the document, merge those three paragraphs together
# (paragraph-mark deletion keeps each line in its own run), and finally
# delete the old "Requisitos" paragraph, leaving the freshly built one
# in its place.

$d = $word.ActiveDocument
$vt = [char]11  # vertical tab == the in-paragraph line break Word stores for <w:br/>

$s1 = "LOB1012 -  Estatística  (Requisito fraco)"
$s2 = "LOQ4095 -  Química Geral Experimental  (Requisito fraco)"
$s3 = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)"

# Locate the existing "Requisitos" list paragraph (the last paragraph,
# containing all three lines) so we can find/replace it later.
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext.IndexOf("LOB1012") -ge 0 -and $ptext.IndexOf("LOQ4095") -ge 0 -and $ptext.IndexOf("LOQ4098") -ge 0) {
        $targetIndex = $i
    }
}

# Step 1: append three placeholder paragraphs at the end of the
# document, in the desired final line order (s2, s3, s1).
$endOfDoc = $d.Content.End
$insA = $d.Range($endOfDoc, $endOfDoc)
$insA.InsertParagraphAfter()
$pA = $d.Paragraphs.Item($d.Paragraphs.Count)
$pA.Range.Text = "@@REQ_PLACEHOLDER_A@@"

$endOfDoc2 = $d.Content.End
$insB = $d.Range($endOfDoc2, $endOfDoc2)
$insB.InsertParagraphAfter()
$pB = $d.Paragraphs.Item($d.Paragraphs.Count)
$pB.Range.Text = "@@REQ_PLACEHOLDER_B@@"

$endOfDoc3 = $d.Content.End
$insC = $d.Range($endOfDoc3, $endOfDoc3)
$insC.InsertParagraphAfter()
$pC = $d.Paragraphs.Item($d.Paragraphs.Count)
$pC.Range.Text = "@@REQ_PLACEHOLDER_C@@"

$firstNewIndex = $targetIndex + 1

# Step 2: swap in the real text (plus trailing line-break char) via
# Find/Replace, which -- unlike a direct Range.Text assignment -- does
# not tag the run with a spurious xml:space="preserve".
$d.Content.Find.Execute("@@REQ_PLACEHOLDER_A@@", $true, $false, $false, $false, $false, $true, 1, $false, ($s2 + $vt), 2) | Out-Null
$d.Content.Find.Execute("@@REQ_PLACEHOLDER_B@@", $true, $false, $false, $false, $false, $true, 1, $false, ($s3 + $vt), 2) | Out-Null
$d.Content.Find.Execute("@@REQ_PLACEHOLDER_C@@", $true, $false, $false, $false, $false, $true, 1, $false, ($s1 + $vt), 2) | Out-Null

# Step 3: merge the three new paragraphs into a single paragraph by
# deleting the paragraph marks between them. Each line keeps its own
# run because the merge happens at the paragraph-mark level rather
# than via a second edit inside an already-edited run.
$firstNew = $d.Paragraphs.Item($firstNewIndex)
$mark1 = $d.Range($firstNew.Range.End - 1, $firstNew.Range.End)
$mark1.Select()
$word.Selection.Delete()

$firstNewAfterMerge = $d.Paragraphs.Item($firstNewIndex)
$mark2 = $d.Range($firstNewAfterMerge.Range.End - 1, $firstNewAfterMerge.Range.End)
$mark2.Select()
$word.Selection.Delete()

# Step 4: remove the original "Requisitos" list paragraph (old line
# order) entirely, including its own trailing paragraph mark. The
# newly built, correctly-ordered paragraph shifts up into its place.
$oldPara = $d.Paragraphs.Item($targetIndex)
$oldRange = $d.Range($oldPara.Range.Start, $oldPara.Range.End)
$oldRange.Select()
$word.Selection.Delete()
